# Update the "取得日時" (acquisition datetime) column (A) for data rows 2-14
# on the "ランサーズ" sheet to reflect the new run timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-26 01:56:47"

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
